# Updates the cryptos worksheet with refreshed price / volume(1h) data.
# For cells whose new value parses as a plain number (e.g. "0.999"),
# force NumberFormat to Text ("@") first so Excel keeps the value as a
# text string (matching the source data, which stores prices as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.333.54'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.514.69'
$ws.Range("E3").Value = '  -2.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.77'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("E6").Value = '  -2.85%  '
$ws.Range("D7").Value = '3.502.72'
$ws.Range("E7").Value = '  -2.23%  '
$ws.Range("E8").Value = '  -2.83%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.187'
$ws.Range("E10").Value = '  +3.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.651'
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.29'
$ws.Range("E12").Value = '  -2.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000302'
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.47'
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").Value = '4.074.34'
$ws.Range("E15").Value = '  -2.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.36'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").Value = '69.217.03'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").Value = '3.507.41'
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '546.70'
$ws.Range("E21").Value = '  +15.17%  '
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.49'
$ws.Range("E23").Value = '  -4.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.95'
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.19'
$ws.Range("E26").Value = '  -1.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.36'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("E28").Value = '  -1.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.15'
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.76'
$ws.Range("E30").Value = '  -1.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.28'
$ws.Range("E31").Value = '  -4.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.69'
$ws.Range("E32").Value = '  +3.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '64.80'
$ws.Range("E34").Value = '  -3.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '546.00'
$ws.Range("E35").Value = '  -7.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.407'
$ws.Range("E36").Value = '  +3.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.20'
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("E38").Value = '  +8.04%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("D40").Value = '0.0₃0767'
$ws.Range("E40").Value = '  -4.25%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.38'
$ws.Range("E41").Value = '  -2.34%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.09'
$ws.Range("E42").Value = '  -3.15%  '
$ws.Range("E43").Value = '  -2.80%  '
$ws.Range("D44").Value = '3.309.36'
$ws.Range("E44").Value = '  +2.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.00'
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0446'
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("E47").Value = '  +3.29%  '
$ws.Range("E48").Value = '  -2.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.95'
$ws.Range("E49").Value = '  -5.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.35'
$ws.Range("E51").Value = '  +2.73%  '
